# edit.ps1
# Applies the "data updates from latest ic results, plus changes in boundary conditions" commit:
#  1. Appends 47 new analytical-result rows (74-120) to Sheet1
#  2. Extends the Sheet1 VLOOKUP range used on the "sampling" sheet from $E$73 to $E$400
#  3. Adds/refreshes the NO3-/SO4-2 (columns M/N) computed columns on "sampling" for the
#     rows whose sample id now resolves against the newly-added Sheet1 rows

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sampling = $wb.Worksheets.Item("sampling")

# ----- 1) Copy down the number formats from the last existing row (73) onto the new rows -----
$sheet1.Range("A73:E73").Copy()
$sheet1.Range("A74:E120").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ----- 2) New Sheet1 data rows 74-120 -----
# columns: row, A (sample label), C, D, E
$sheet1NewData = @(
    ,@(74, "008-F1-62", 14.226351876830869, 222.31040627747802, 0.13316378071010068)
    ,@(75, "009-F2-62", 14.191341775350649, 215.72683116526787, 0.14608727899638632)
    ,@(76, "010-F3-62", 14.205443592326224, 230.03500705713708, 0.14186298530970723)
    ,@(77, "011-F4-62", 14.199060334152474, 230.77163431382317, 0.12255515900393749)
    ,@(78, "012-F1-59", 14.83243417635699, 223.11333641621096, 0.13357809354702721)
    ,@(79, "013-F2-59", 14.731503774128422, 219.05562138639962, 0.10422868399678681)
    ,@(80, "014-F3-59", 14.66241234508437, 232.0577134677041, 0.1153421798310512)
    ,@(81, "015-F4-59", 15.438007115511359, 232.33512183528407, 0.1072978023933663)
    ,@(82, "016-F1-57", 14.206518322564211, 222.59276834506252, 0.12147748257612119)
    ,@(83, "017-F2-57", 14.396573646771763, 221.96672930627574, 0.1249590480826368)
    ,@(84, "018-F3-57", 14.152129733484934, 232.27201176317521, 0.11791254751119029)
    ,@(85, "019-F1-56", 14.547797352410246, 228.83649911025572, 0.22233926079185029)
    ,@(86, "020-F2-56", 14.19029960478913, 224.55323864462892, 0.13788656173788003)
    ,@(87, "021-F3-56", 14.371108063794351, 234.83967076096826, 0.1508905558534247)
    ,@(88, "022-F4-56", 14.291224787671471, 234.16530063745256, 0.1571831809419327)
    ,@(89, "023-F1-55", 15.713144691566701, 228.95397140920653, 0.1814288516507708)
    ,@(90, "024-F2-55", 17.439259705141964, 229.49386012837763, 0.12562215141813121)
    ,@(91, "025-F3-55", 17.878681867224639, 234.99070753853155, 0.13772086442261722)
    ,@(92, "026-F4-55", 18.250014864258027, 234.39657297524838, 0.14111745151918231)
    ,@(93, "027-F1-53", 14.416437878733555, 231.10276627067952, 0.17729294465374082)
    ,@(94, "028-F2-53", 14.425360436677311, 230.46329757256197, 0.12893741850001922)
    ,@(95, "029-F3-53", 14.339845523227581, 235.53384216024318, 0.1547822248713408)
    ,@(96, "030-F4-53", 14.410413501381923, 235.24889049085141, 0.14385097245959519)
    ,@(97, "031-F1-49", 13.214068277703049, 249.13725106757619, 0.12379857720623)
    ,@(98, "032-F2-49", 13.220682766992441, 247.55618257314066, 0.12960042199727997)
    ,@(99, "033-F3-49", 13.223159121887248, 248.26760189999709, 0.16877092760610471)
    ,@(100, "034-F4-49", 13.185329208199558, 249.29543370935659, 0.12214067327007)
    ,@(101, "035-F1-51", 14.392600776679016, 239.33028294095988, 0.12844015495726682)
    ,@(102, "036-F2-51", 14.47264281585552, 237.29876576298457, 0.1131861949266012)
    ,@(103, "037-F3-51", 14.472610252618576, 240.31421891618737, 0.12238936549037671)
    ,@(104, "038-F4-51", 14.423634543541228, 241.16271620948774, 0.10945426724932469)
    ,@(105, "039-F1-52", 14.47938339447602, 237.11349162551934, 0.14923435228718068)
    ,@(106, "040-F2-52", 14.571469794604063, 233.15378503825872, 0.14459642903148628)
    ,@(107, "041-F3-52", 14.485179622638841, 239.51802476251896, 0.12479326964880919)
    ,@(108, "042-F4-52", 14.445940783331913, 237.06556263227267, 0.14492773630355749)
    ,@(109, "043-F1-48", 13.569395034472938, 248.07309348138384, 0.12255515900393749)
    ,@(110, "044-F2-48", 13.817981964522671, 248.74082507842834, 0.1168346818067228)
    ,@(111, "045-F3-48", 13.741399244920672, 249.8373891970503, 0.1459216331602127)
    ,@(112, "046-F4-48", 13.635918114559018, 250.92342462093254, 0.14376814265163029)
    ,@(113, "047-F1-45", 12.966612013782999, 249.21736979897349, 0.13763801537499271)
    ,@(114, "048-F2-45", 13.024943281946751, 249.57884596854538, 0.1373894666721468)
    ,@(115, "049-F3-45", 12.970294429048135, 250.22138710608178, 0.1389635689375367)
    ,@(116, "050-F4-45", 13.116770848998399, 252.46004013629351, 0.21110557619715753)
    ,@(117, "051-F1-44", 19.14984585958593, 249.16154779304591, 0.1398748484024688)
    ,@(118, "052-F2-44", 53.744562692455446, 252.83993212089996, 0.13722376623693922)
    ,@(119, "053-F3-44", 16.289495765805054, 250.54300260120752, 0.11326912067286871)
    ,@(120, "054-F4-44", 15.197629372884938, 251.33168362264513, 0.18061977447124519)
)

foreach ($row in $sheet1NewData) {
    $r = $row[0]
    $sheet1.Cells.Item($r, 1).Value = $row[1]
    $sheet1.Cells.Item($r, 2).Formula = "=MID(A$r,5,5)"
    $sheet1.Cells.Item($r, 3).Value = $row[2]
    $sheet1.Cells.Item($r, 4).Value = $row[3]
    $sheet1.Cells.Item($r, 5).Value = $row[4]
}

# ----- 3) sampling sheet: (re)write columns M (NO3-) and N (SO4-2) -----
# Row 165 already had these formulas (pointing at the old $E$73 range); the rest are brand new
# cells made reachable now that their sample id exists in Sheet1.
$samplingMNRows = @(165,174,175,176,177,178,179,180,181,190,191,192,193,194,195,196,197,202,203,204,205,206,207,208,209,210,211,212,213,218,219,220,221,222,223,224,225,226,227,228,234,235,236,237,246,247,248,249)

foreach ($r in $samplingMNRows) {
    $sampling.Cells.Item($r, 13).Formula = "=VLOOKUP(C$r,Sheet1!`$B`$3:`$E`$400,4,FALSE)/62*1000"
    $sampling.Cells.Item($r, 14).Formula = "=VLOOKUP(C$r,Sheet1!`$B`$3:`$E`$400,3,FALSE)/96*1000"
}

# ----- 4) sheet view / selection cosmetics matching the saved workbook state -----
$samplingWindow = $excel.Windows.Item(1)
$sheet1.Activate()
$sheet1.Range("I115").Select()
$excel.ActiveWindow.ScrollRow = 62

$sampling.Activate()
$sampling.Range("M166:N173").Select()
$excel.ActiveWindow.ScrollColumn = 2

Write-Host "Applied no2_analyticall_results updates"